# Fixes the label/value misalignment in the LOT2007 "Bioquimica I" syllabus sheet.
# The sheet was missing a row for "Docentes responsaveis:" value, which had pushed
# every value below it up by one row relative to its label. This inserts the missing
# row and refreshes the objectives / summary / syllabus / evaluation / bibliography text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 - this shifts old rows 13:24 down to 14:25 and keeps all
# their row heights / A-column labels correctly aligned with the row below them.
$ws.Rows("13:13").Insert()

# New row 13 holds the (previously missing) value for "Docentes responsaveis:" (row 12),
# so it has no column-A label of its own - copy the B:C number format down, then clear A.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '427823 - Adriane Maria Ferreira Milagres'
$ws.Range("C13").Value = '427823 - Adriane Maria Ferreira Milagres'

# Refresh the text that now needs to change in place (same label, new content).
# Row 10: Objetivos: (new long-form text)
$ws.Range("B10").Value = 'Promover aos participantes do curso conhecimentos de bioquímica  abrangendo a organização estrutural e molecular da célulaCompreender a importância dos compostos orgânicos no metabolismo celularUtilizar todos os conhecimentos como pré-requisito para as disciplinas do curso de engenharia Bioquímica'
$ws.Range("C10").Value = 'Promover aos participantes do curso conhecimentos de bioquímica  abrangendo a organização estrutural e molecular da célulaCompreender a importância dos compostos orgânicos no metabolismo celularUtilizar todos os conhecimentos como pré-requisito para as disciplinas do curso de engenharia Bioquímica'

# Row 14: Programa resumido: (new short program outline)
$ws.Range("B14").Value = '01Química ácido-base/Tampões02Aminoácidos03Proteínas:Estrutura primária04Proteínas:Estrutura tridimensional05Função das proteínas06Enzimas: catálise enzimática07Cinética enzimática, inibição e regulação08Carboidratos09Lipídeos10Membranas Biológicas11Nucleotídeos e ácidos nucleicos'
$ws.Range("C14").Value = '01Química ácido-base/Tampões02Aminoácidos03Proteínas:Estrutura primária04Proteínas:Estrutura tridimensional05Função das proteínas06Enzimas: catálise enzimática07Cinética enzimática, inibição e regulação08Carboidratos09Lipídeos10Membranas Biológicas11Nucleotídeos e ácidos nucleicos'

# Row 16: Programa: (new full program description)
$ws.Range("B16").Value = '01Química ácido-base/Tampões : Constante de dissociação, Curvas de titulação , capacidade tamponante02Aminoácidos:Estrutura dos aminoácidos, classificação e características, Nomenclatura, propriedades ácido-base, estereoquímica,aminoácidos incomuns03Proteínas:Estrutura primáriaPurificação de proteínas, solubilidade, cormatografia, eletroforese. Sequenciamento de proteínas.04Proteínas:Estrutura tridimensionalEstrutura secundária, terciária, quaternária. Dobramento e estabilidade das proteínas.05Função das proteínasMioglobina, hemoglobina, anticorpos06Enzimas: catálise enzimáticaNomenclatura das enzimas, especificidade dos substratos, co-fatores e coenzimas, Energia de ativação e coordenada de reação. Curvas de progresso. Efeito da temperatura e pH sobre a velocidade das reações enzimáticas07Cinética enzimática, inibição e regulação Efeito da concentração de substrato na velocidade das reações enzimáticas. Efeito de inibidores. Inibições reversíveis.  Modelos de inibição competitiva, não competitiva e acompetitiva simples.08CarboidratosMonossacarídeos:classificação, configuração e conformação. Dissacarídeos, Polissacarídeos estruturais:celulose e quitina, Polissacarídeos de reserva:amido e glicogênio, glicosaminoglicanos, Glicoproteínas:oligossacarídeos, paredes celulares bacterianas.09LipídeosClassificação: ácidos graxos, trioacilglicerol, glicerofosfolipídeos, esfingolipídeos, esteróides.10Membranas BiológicasProteínas de membrana:integrais e periféricas, modelo do mosaico fluido, assimetria dos lipídeos, Transporte através da membrana: termodinâmica do transporte,  transporte passivo e ativo.11Nucleotídeos e ácidos nucleicos:Estrutura e função dos nucleotídeos, Estrutura dos ácidos nucleicos, sequenciamento de ácidos nucleicos, endonucleases de restrição, Bibliotecas genômicas,Amplificação do DNA pela reação em cadeia da polimerase.'
$ws.Range("C16").Value = '01Química ácido-base/Tampões : Constante de dissociação, Curvas de titulação , capacidade tamponante02Aminoácidos:Estrutura dos aminoácidos, classificação e características, Nomenclatura, propriedades ácido-base, estereoquímica,aminoácidos incomuns03Proteínas:Estrutura primáriaPurificação de proteínas, solubilidade, cormatografia, eletroforese. Sequenciamento de proteínas.04Proteínas:Estrutura tridimensionalEstrutura secundária, terciária, quaternária. Dobramento e estabilidade das proteínas.05Função das proteínasMioglobina, hemoglobina, anticorpos06Enzimas: catálise enzimáticaNomenclatura das enzimas, especificidade dos substratos, co-fatores e coenzimas, Energia de ativação e coordenada de reação. Curvas de progresso. Efeito da temperatura e pH sobre a velocidade das reações enzimáticas07Cinética enzimática, inibição e regulação Efeito da concentração de substrato na velocidade das reações enzimáticas. Efeito de inibidores. Inibições reversíveis.  Modelos de inibição competitiva, não competitiva e acompetitiva simples.08CarboidratosMonossacarídeos:classificação, configuração e conformação. Dissacarídeos, Polissacarídeos estruturais:celulose e quitina, Polissacarídeos de reserva:amido e glicogênio, glicosaminoglicanos, Glicoproteínas:oligossacarídeos, paredes celulares bacterianas.09LipídeosClassificação: ácidos graxos, trioacilglicerol, glicerofosfolipídeos, esfingolipídeos, esteróides.10Membranas BiológicasProteínas de membrana:integrais e periféricas, modelo do mosaico fluido, assimetria dos lipídeos, Transporte através da membrana: termodinâmica do transporte,  transporte passivo e ativo.11Nucleotídeos e ácidos nucleicos:Estrutura e função dos nucleotídeos, Estrutura dos ácidos nucleicos, sequenciamento de ácidos nucleicos, endonucleases de restrição, Bibliotecas genômicas,Amplificação do DNA pela reação em cadeia da polimerase.'

# Row 19: Metodo: (avaliacao criteria text moved up)
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

# Row 20: Criterio: (nota final formula moved up)
$ws.Range("B20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3'
$ws.Range("C20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3'

# Row 21: Norma de recuperacao: (recovery formula moved up)
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'

# Row 22: Bibliografia: (new bibliography text)
$ws.Range("B22").Value = '1. M. Cox, Michael; Nelson, David L.Princípios de Bioquímica de Lehninger - Editora Artmed 6ª Ed. 2014
2. Voet, D., Voet, J. G., Pratt, C.W. Fundamentos de Bioquímica:a vida em nivel molecular  Editora Artmed, 2014'
$ws.Range("C22").Value = '1. M. Cox, Michael; Nelson, David L.Princípios de Bioquímica de Lehninger - Editora Artmed 6ª Ed. 2014
2. Voet, D., Voet, J. G., Pratt, C.W. Fundamentos de Bioquímica:a vida em nivel molecular  Editora Artmed, 2014'

# Column A no longer shares its width definition with column B (A stays 30.71,
# B keeps its own 60.71 override) - touch column B so the engine splits the range.
$ws.Columns("B").ColumnWidth = 60.7109375

